$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data permutation: each row's Fecha/Volumen/Precio fields are reshuffled
# (weekly re-aggregation). $data maps target row -> @(D, J, K, L, M, P).
$data = @{
    2 = @(44225, 56, 3000, 3000, 3000, 1000)
    3 = @(44292, 40, 3000, 3000, 3000, 1000)
    4 = @(44242, 95, 2500, 3000, 2737, 912)
    5 = @(44193, 70, 3000, 3000, 3000, 1000)
    6 = @(44557, 104, 2000, 2500, 2260, 753)
    7 = @(44574, 50, 3000, 3000, 3000, 1000)
    8 = @(44291, 45, 3000, 3000, 3000, 1000)
    9 = @(44559, 68, 2000, 2000, 2000, 667)
    10 = @(44223, 80, 2500, 3000, 2781, 927)
    11 = @(44390, 50, 3000, 3000, 3000, 1000)
    12 = @(44166, 45, 2500, 2500, 2500, 833)
    13 = @(44340, 54, 3000, 3000, 3000, 1000)
    14 = @(44536, 125, 2200, 2200, 2200, 733)
    15 = @(44222, 45, 3000, 3000, 3000, 1000)
    16 = @(44627, 78, 3500, 3500, 3500, 1167)
    17 = @(44260, 60, 3500, 3500, 3500, 1167)
    18 = @(44243, 45, 3000, 3000, 3000, 1000)
    19 = @(44187, 65, 3000, 3000, 3000, 1000)
    20 = @(44669, 92, 2500, 3000, 2755, 918)
    21 = @(44389, 81, 2800, 3000, 2889, 963)
    22 = @(44179, 78, 3000, 3000, 3000, 1000)
    23 = @(44224, 67, 3000, 3000, 3000, 1000)
    24 = @(44537, 88, 2000, 2200, 2091, 697)
    25 = @(44165, 68, 3000, 3000, 3000, 1000)
    26 = @(44221, 50, 2500, 2500, 2500, 833)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
    $ws.Range("K$row").Value = $vals[2]
    $ws.Range("L$row").Value = $vals[3]
    $ws.Range("M$row").Value = $vals[4]
    $ws.Range("P$row").Value = $vals[5]
}
